$wb = $excel.ActiveWorkbook

# Sheet1 (weibull)
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.80076969518453
$ws.Range("C2").Value = 0.199420834984694
$ws.Range("B3").Value = 0.300144465918335
$ws.Range("C3").Value = 0.135723151780683

# Sheet2 (lognormal)
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.39272083727114
$ws.Range("C2").Value = 0.265967582835733
$ws.Range("B3").Value = -1.21578962659053
$ws.Range("C3").Value = 0.119222995723348

# Sheet3 (llogis)
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.99064811632014
$ws.Range("C2").Value = 0.0821751582283421
$ws.Range("B3").Value = 0.785020862239596
$ws.Range("C3").Value = 0.106565176953077

# Sheet4 (gompertz)
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.34740021745688
$ws.Range("C2").Value = 0.126118366223994
$ws.Range("B3").Value = 0.00794657992825708
$ws.Range("C3").Value = 0.0201871014722952

# Sheet6 (weibull cov)
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0397686694259927
$ws.Range("B2").Value = -0.0225109028945209
$ws.Range("A3").Value = -0.0225109028945209
$ws.Range("B3").Value = 0.0184207739292822

# Sheet7 (lognormal cov)
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0707387551194827
$ws.Range("B2").Value = -0.0289436015331549
$ws.Range("A3").Value = -0.0289436015331549
$ws.Range("B3").Value = 0.0142141227092495

# Sheet8 (llogis cov)
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00675275662985306
$ws.Range("B2").Value = -0.00158573335259087
$ws.Range("A3").Value = -0.00158573335259087
$ws.Range("B3").Value = 0.0113561369390407

# Sheet9 (gompertz cov)
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0159058422990094
$ws.Range("B2").Value = -0.0017874632243433
$ws.Range("A3").Value = -0.0017874632243433
$ws.Range("B3").Value = 0.000407519065852742
